$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1572.56
$ws.Range("I28").Value = 1817.381
$ws.Range("J28").Value = 287.25
$ws.Range("K28").Value = 1817.381
$ws.Range("L28").Value = 287.25
$ws.Range("M28").Value = -1332.381
$ws.Range("N28").Value = -1257.25

$ws.Range("H38").Value = 1611.2195
$ws.Range("I38").Value = 212.875
$ws.Range("J38").Value = 1950.2122
$ws.Range("K38").Value = 638.625
$ws.Range("L38").Value = 5850.6366
$ws.Range("M38").Value = -266.625
$ws.Range("N38").Value = -6594.6366

$ws.Range("H58").Value = 1600.6666
$ws.Range("J58").Value = 3577.125
$ws.Range("L58").Value = 10731.375
$ws.Range("N58").Value = -11031.375

$ws.Range("H74").Value = 3126.5715
$ws.Range("I74").Value = 2857.2
$ws.Range("K74").Value = 2857.2
$ws.Range("M74").Value = -1921.2

$ws.Range("H77").Value = 3126.5715
$ws.Range("I77").Value = 2857.2
$ws.Range("K77").Value = 14286
$ws.Range("M77").Value = -9606

$ws.Range("H106").Value = 6072.815
$ws.Range("I106").Value = 6628.4346
$ws.Range("K106").Value = 6628.4346
$ws.Range("M106").Value = -5997.4346

$ws.Range("H113").Value = 3702.75
$ws.Range("I113").Value = 1752.5
$ws.Range("J113").Value = 5653
$ws.Range("K113").Value = 1752.5
$ws.Range("L113").Value = 5653
$ws.Range("M113").Value = 1501.5
$ws.Range("N113").Value = -12161

$ws.Range("H132").Value = 6950289
$ws.Range("I132").Value = 10422180
$ws.Range("J132").Value = 6505.75
$ws.Range("K132").Value = 31266540
$ws.Range("L132").Value = 19517.25
$ws.Range("M132").Value = -31264010
$ws.Range("N132").Value = -24577.25

$ws.Range("H138").Value = 821974.1
$ws.Range("I138").Value = 1370.3529
$ws.Range("J138").Value = 1209481.5
$ws.Range("K138").Value = 4111.0587
$ws.Range("L138").Value = 3628444.5
$ws.Range("M138").Value = 1028.9413
$ws.Range("N138").Value = -3638724.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5097.82
$ws.Range("I32").Value = 4088.4465
$ws.Range("K32").Value = 4088.4465
$ws.Range("M32").Value = -3801.4465

$ws.Range("H61").Value = 854.1667
$ws.Range("I61").Value = 709.1875
$ws.Range("K61").Value = 709.1875
$ws.Range("M61").Value = -497.1875

$ws.Range("H74").Value = 1731.4073
$ws.Range("J74").Value = 2873.5
$ws.Range("L74").Value = 2873.5
$ws.Range("N74").Value = -4621.5

$ws.Range("H77").Value = 1731.4073
$ws.Range("J77").Value = 2873.5
$ws.Range("L77").Value = 14367.5
$ws.Range("N77").Value = -23103.5

$ws.Range("H110").Value = 4406.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 4406.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 4406.5
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -8496.5

$ws.Range("H132").Value = 3123.4138
$ws.Range("I132").Value = 2575.182
$ws.Range("J132").Value = 4846.4287
$ws.Range("K132").Value = 7725.545999999999
$ws.Range("L132").Value = 14539.2861
$ws.Range("M132").Value = -5195.545999999999
$ws.Range("N132").Value = -19599.2861

$ws.Range("H136").Value = 854.1667
$ws.Range("I136").Value = 709.1875
$ws.Range("K136").Value = 2127.5625
$ws.Range("M136").Value = 422.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 510.9375
$ws.Range("J80").Value = 607.6923
$ws.Range("L80").Value = 607.6923
$ws.Range("N80").Value = -2603.6923

$ws.Range("H83").Value = 510.9375
$ws.Range("J83").Value = 607.6923
$ws.Range("L83").Value = 3038.4615
$ws.Range("N83").Value = -13022.4615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 213.25
$ws.Range("I22").Value = 184.33333
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 184.33333
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 165.66667
$ws.Range("N22").Value = -1000

$ws.Range("H58").Value = 1776.7368
$ws.Range("I58").Value = 1424.5454
$ws.Range("K58").Value = 1424.5454
$ws.Range("M58").Value = -1221.5454

$ws.Range("H62").Value = 8697891
$ws.Range("J62").Value = 100001096
$ws.Range("L62").Value = 100001096
$ws.Range("N62").Value = -100002344

$ws.Range("H65").Value = 8697891
$ws.Range("J65").Value = 100001096
$ws.Range("L65").Value = 500005480
$ws.Range("N65").Value = -500011720

$ws.Range("H105").Value = 534.8333
$ws.Range("I105").Value = 481.8
$ws.Range("K105").Value = 481.8
$ws.Range("M105").Value = 1265.2

$ws.Range("H132").Value = 1705.2368
$ws.Range("I132").Value = 1314.5385
$ws.Range("K132").Value = 3943.6155
$ws.Range("M132").Value = -1413.6155

$ws.Range("H136").Value = 1776.7368
$ws.Range("I136").Value = 1424.5454
$ws.Range("K136").Value = 4273.6362
$ws.Range("M136").Value = -1723.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 47619172
$ws.Range("I12").Value = 166666850
$ws.Range("J12").Value = 100.86667
$ws.Range("K12").Value = 500000550
$ws.Range("L12").Value = 302.60001
$ws.Range("M12").Value = -500000377
$ws.Range("N12").Value = -648.60001

$ws.Range("H39").Value = 2701.9524
$ws.Range("I39").Value = 1950
$ws.Range("J39").Value = 2781.1052
$ws.Range("K39").Value = 5850
$ws.Range("L39").Value = 8343.3156
$ws.Range("M39").Value = -5556
$ws.Range("N39").Value = -8931.3156

$ws.Range("H55").Value = 2349.5
$ws.Range("I55").Value = 875
$ws.Range("J55").Value = 3332.5
$ws.Range("K55").Value = 2625
$ws.Range("L55").Value = 9997.5
$ws.Range("M55").Value = -2448
$ws.Range("N55").Value = -10351.5

$ws.Range("H76").Value = 6001
$ws.Range("I76").Value = 5506.5
$ws.Range("J76").Value = 6090.909
$ws.Range("K76").Value = 16519.5
$ws.Range("L76").Value = 18272.727
$ws.Range("M76").Value = -16136.5
$ws.Range("N76").Value = -19038.727

$ws.Range("H79").Value = 6001
$ws.Range("I79").Value = 5506.5
$ws.Range("J79").Value = 6090.909
$ws.Range("K79").Value = 16519.5
$ws.Range("L79").Value = 18272.727
$ws.Range("M79").Value = -15193.5
$ws.Range("N79").Value = -20924.727

$ws.Range("H117").Value = 1249.4166
$ws.Range("J117").Value = 1457.3334
$ws.Range("L117").Value = 4372.0002
$ws.Range("N117").Value = -11256.0002

$ws.Range("H121").Value = 875
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -5620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3141.6667
$ws.Range("I80").Value = 2431.4167
$ws.Range("J80").Value = 4088.6667
$ws.Range("K80").Value = 2431.4167
$ws.Range("L80").Value = 4088.6667
$ws.Range("M80").Value = -1433.4167
$ws.Range("N80").Value = -6084.6667

$ws.Range("H83").Value = 3141.6667
$ws.Range("I83").Value = 2431.4167
$ws.Range("J83").Value = 4088.6667
$ws.Range("K83").Value = 12157.0835
$ws.Range("L83").Value = 20443.3335
$ws.Range("M83").Value = -7165.083500000001
$ws.Range("N83").Value = -30427.3335

$ws.Range("H132").Value = 3888.2
$ws.Range("I132").Value = 3841.1
$ws.Range("K132").Value = 11523.3
$ws.Range("M132").Value = -8993.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2725.8667
$ws.Range("I132").Value = 1600
$ws.Range("J132").Value = 2899.077
$ws.Range("K132").Value = 4800
$ws.Range("L132").Value = 8697.231
$ws.Range("M132").Value = -2270
$ws.Range("N132").Value = -13757.231

$ws.Range("H136").Value = 2172
$ws.Range("I136").Value = 1534.6666
$ws.Range("J136").Value = 2650
$ws.Range("K136").Value = 4603.9998
$ws.Range("L136").Value = 7950
$ws.Range("M136").Value = -2053.9998
$ws.Range("N136").Value = -13050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
